# New staging template field: "ActualDate" is added as the first data
# column on the header row (row 2), pushing the existing headers
# (ActualLabel .. StrategicElementBusinessKey) one column to the right
# (B2:U2). Row 1 (the "for internal use only" banner cell, A1) is left
# untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastCol = 20  # T2 is currently the last populated header cell in row 2

# Shift existing row-2 header values one column to the right, starting
# from the rightmost column so we never clobber a value before it has
# been copied.
for ($c = $lastCol; $c -ge 1; $c--) {
    $srcCell = $ws.Cells.Item(2, $c)
    $dstCell = $ws.Cells.Item(2, $c + 1)
    $dstCell.Value = $srcCell.Value2
}

# New first header column.
$ws.Range("A2").Value = "ActualDate"
